# Thu Jan 19 11:35:17 UTC 2023 GitHub Actions refresh of the coin table.
# Column D = Price, Column E = Volume(1h) (a percentage rendered as literal text).
# Both columns are stored as plain text in this sheet (e.g. "0.1000", "-3.63%"), so
# the cells are forced to Text format before the write to stop Excel's automatic
# number/percentage inference from collapsing the exact published formatting
# (trailing zeros, "%" suffix, sign, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "290.06" },
    @{ Cell = "E2"; Value = "-3.63%" },
    @{ Cell = "D3"; Value = "30.44" },
    @{ Cell = "E3"; Value = "-5.59%" },
    @{ Cell = "D4"; Value = "4.942" },
    @{ Cell = "E4"; Value = "-1.01%" },
    @{ Cell = "D5"; Value = "0.07207" },
    @{ Cell = "E5"; Value = "-5.54%" },
    @{ Cell = "D6"; Value = "1.826" },
    @{ Cell = "E6"; Value = "-6.28%" },
    @{ Cell = "D7"; Value = "7.678" },
    @{ Cell = "E7"; Value = "-2.45%" },
    @{ Cell = "D8"; Value = "3.759" },
    @{ Cell = "E8"; Value = "-0.65%" },
    @{ Cell = "D9"; Value = "0.8973" },
    @{ Cell = "E9"; Value = "-2.29%" },
    @{ Cell = "D10"; Value = "0.1657" },
    @{ Cell = "E10"; Value = "-6.52%" },
    @{ Cell = "D11"; Value = "0.07717" },
    @{ Cell = "E11"; Value = "-1.39%" },
    @{ Cell = "D12"; Value = "0.07958" },
    @{ Cell = "E12"; Value = "-6.48%" },
    @{ Cell = "D13"; Value = "0.03038" },
    @{ Cell = "E13"; Value = "-3.96%" },
    @{ Cell = "D14"; Value = "0.1000" },
    @{ Cell = "E14"; Value = "0.04%" },
    @{ Cell = "D15"; Value = "0.001497" },
    @{ Cell = "E15"; Value = "-1.22%" },
    @{ Cell = "D16"; Value = "0.005878" },
    @{ Cell = "E16"; Value = "-0.13%" },
    @{ Cell = "D18"; Value = "3.463" },
    @{ Cell = "E18"; Value = "0.04%" },
    @{ Cell = "D19"; Value = "2.079" },
    @{ Cell = "E19"; Value = "-3.32%" },
    @{ Cell = "D20"; Value = "0.3320" },
    @{ Cell = "E20"; Value = "-0.68%" },
    @{ Cell = "E21"; Value = "-1.53%" },
    @{ Cell = "D22"; Value = "4.047" },
    @{ Cell = "E22"; Value = "-5.15%" },
    @{ Cell = "D23"; Value = "0.2391" },
    @{ Cell = "E23"; Value = "20.10%" },
    @{ Cell = "D24"; Value = "0.04518" },
    @{ Cell = "E24"; Value = "0.18%" },
    @{ Cell = "D25"; Value = "0.001215" },
    @{ Cell = "E25"; Value = "-0.38%" },
    @{ Cell = "D26"; Value = "0.004616" },
    @{ Cell = "E26"; Value = "5.22%" },
    @{ Cell = "E27"; Value = "4.11%" },
    @{ Cell = "D39"; Value = "0.01563" },
    @{ Cell = "E39"; Value = "-8.08%" },
    @{ Cell = "D40"; Value = "0.04354" },
    @{ Cell = "E40"; Value = "-6.86%" },
    @{ Cell = "D41"; Value = "0.007318" },
    @{ Cell = "E41"; Value = "-2.47%" },
    @{ Cell = "D42"; Value = "0.009817" },
    @{ Cell = "E43"; Value = "-3.51%" },
    @{ Cell = "D44"; Value = "0.002016" },
    @{ Cell = "E44"; Value = "-13.54%" },
    @{ Cell = "D45"; Value = "0.009513" },
    @{ Cell = "E45"; Value = "-8.86%" },
    @{ Cell = "D46"; Value = "0.00006016" },
    @{ Cell = "E46"; Value = "-3.81%" },
    @{ Cell = "E47"; Value = "0.13%" },
    @{ Cell = "D48"; Value = "2.255" },
    @{ Cell = "E48"; Value = "174.84%" },
    @{ Cell = "D50"; Value = "0.00002104" },
    @{ Cell = "E50"; Value = "0.13%" },
    @{ Cell = "D51"; Value = "0.0002004" },
    @{ Cell = "E51"; Value = "0.13%" }
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    # Quote-prefix / Text format so the numeric-looking string round-trips as
    # inline text instead of being reinterpreted as a Number/Percentage.
    $r.NumberFormat = "@"
    $r.Value = $u.Value
}
